$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.038.17"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.89%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.818.03"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.09%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.018"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +1.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.97"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.009"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4295"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3679"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07252"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.207.42"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +25.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8673"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.23"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.414"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.617"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06973"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.23"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.015"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008933"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.24%  "
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.26"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.104.29"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.178"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.408.92"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +21.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.99"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.45"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.883"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.44%  "
$ws.Range("E27").Value = "  +0.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.219"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.898"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "114.57"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08962"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.190"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7531"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.427"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.808"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.009"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.136"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05214"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01928"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5117"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.750"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1649"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.481"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.350"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "106.94"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.45"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.010"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4576"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.880"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.646"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06228"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.13%  "
